$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B2").Value = 0.1637010676156584
$summary.Range("C2").Value = 0.05622489959839357
$summary.Range("E2").Value = 0.1064638783269962
$summary.Range("F2").Value = 0.2295081967213115
$summary.Range("G2").Value = 0.6076794657762938
$summary.Range("H2").Value = 0.7920679507758159
$summary.Range("J2").Value = 470
$summary.Range("K2").Value = 64

# --- Classification Report sheet ---
$cr = $wb.Worksheets.Item("Classification Report")
$cr.Range("C2").Value = 0.1198501872659176
$cr.Range("D2").Value = 0.2140468227424749

$cr.Range("B3").Value = 0.05622489959839357
$cr.Range("D3").Value = 0.1064638783269962

$cr.Range("B4").Value = 0.1637010676156584
$cr.Range("C4").Value = 0.1637010676156584
$cr.Range("D4").Value = 0.1637010676156584
$cr.Range("E4").Value = 0.1637010676156584

$cr.Range("B5").Value = 0.5281124497991968
$cr.Range("C5").Value = 0.5599250936329588
$cr.Range("D5").Value = 0.1602553505347356

$cr.Range("B6").Value = 0.9529791764924467
$cr.Range("C6").Value = 0.1637010676156584
$cr.Range("D6").Value = 0.2086868183943728

# --- Confusion Matrix sheet ---
$cm = $wb.Worksheets.Item("Confusion Matrix")
$cm.Range("B2").Value = 64
$cm.Range("C2").Value = 470
